$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("8:8").Insert()
$ws.Range("A8").Value = "被"
$ws.Range("B8").Value = "passive"
$ws.Range("C8").Value = "虛詞"
